# Update "想去人数" (number of people interested) figures to the latest
# scraped counts, as output by the gh-pages generator at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 101
$ws.Range("F5").Value  = 1707
$ws.Range("F6").Value  = 3287
$ws.Range("F7").Value  = 931
$ws.Range("F8").Value  = 2116
$ws.Range("F9").Value  = 2039
$ws.Range("F10").Value = 1049
$ws.Range("F11").Value = 564
$ws.Range("F14").Value = 358
$ws.Range("F18").Value = 116
$ws.Range("F19").Value = 1489
$ws.Range("F20").Value = 560
$ws.Range("F21").Value = 662
$ws.Range("F22").Value = 546
$ws.Range("F23").Value = 11907
$ws.Range("F24").Value = 11932
$ws.Range("F25").Value = 874
$ws.Range("F26").Value = 672
$ws.Range("F27").Value = 57
$ws.Range("F28").Value = 1869
$ws.Range("F29").Value = 170
$ws.Range("F30").Value = 489

# --- Sheet "演出" ---------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 2

# --- Sheet "本地生活" ------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 67

# --- Sheet "全部类型" ------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 67
$ws.Range("F6").Value  = 101
$ws.Range("F7").Value  = 1707
$ws.Range("F8").Value  = 3287
$ws.Range("F9").Value  = 931
$ws.Range("F10").Value = 2116
$ws.Range("F11").Value = 2039
$ws.Range("F12").Value = 1049
$ws.Range("F13").Value = 564
$ws.Range("F16").Value = 358
$ws.Range("F22").Value = 116
$ws.Range("F23").Value = 1489
$ws.Range("F24").Value = 560
$ws.Range("F25").Value = 662
$ws.Range("F26").Value = 546
$ws.Range("F27").Value = 11907
$ws.Range("F28").Value = 11932
$ws.Range("F29").Value = 874
$ws.Range("F30").Value = 672
$ws.Range("F31").Value = 57
$ws.Range("F32").Value = 1869
$ws.Range("F35").Value = 170
$ws.Range("F36").Value = 489
$ws.Range("F37").Value = 2
